$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Edit 1: split the "Para la refrigeracion del edificio..." sentence so
# that " cuerpo principal del " is inserted as its own run between
# "Para la refrigeracion del" and "edificio hemos optado...".
# ---------------------------------------------------------------------
$p2 = $d.Paragraphs.Item(2)
$p2Start = $p2.Range.Start

# Insert the new wording right before "edificio" (offset 26 from the
# paragraph start: "Para la refrigeración del" is 25 chars + 1 space).
$insertionPoint = $d.Range($p2Start + 26, $p2Start + 26)
$insertionPoint.InsertBefore("cuerpo principal del ")

# Force the text between "del" and "edificio" (now " cuerpo principal
# del ") to become its own run by bracketing it with a temporary
# bookmark and immediately deleting the bookmark again; this leaves the
# run split in place without adding any bookmark or formatting marks.
$middle = $d.Range($p2Start + 25, $p2Start + 47)
$d.Bookmarks.Add("tmp_run_split", $middle)
$d.Bookmarks("tmp_run_split").Delete()

# ---------------------------------------------------------------------
# Edit 2: add a new paragraph after the servers paragraph describing the
# air conditioning for the expedition/production building, and move the
# _GoBack bookmark so it still sits at the very end of the document.
# ---------------------------------------------------------------------
$p3 = $d.Paragraphs.Item(3)
$p3.Range.InsertParagraphAfter()

$p4 = $d.Paragraphs.Item(4)
$p4.Range.InsertBefore("Además habrá dos aire acondicionado de 18 mil frigorías en el  edificio de expedición/producción, que funcionarán en la noche, horario en que trabajan las rotativas.")

# Adding a bookmark exactly at the current end-of-document position is
# not persisted reliably, so temporarily append a placeholder character
# after the new text, anchor the bookmark just before it, and then
# remove the placeholder again (the point bookmark stays in place).
$p4 = $d.Paragraphs.Item(4)
$tailPoint = $d.Range($p4.Range.End - 1, $p4.Range.End - 1)
$tailPoint.InsertAfter("Z")

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$p4 = $d.Paragraphs.Item(4)
$bookmarkPoint = $d.Range($p4.Range.End - 2, $p4.Range.End - 2)
$d.Bookmarks.Add("_GoBack", $bookmarkPoint)

$p4 = $d.Paragraphs.Item(4)
$placeholder = $d.Range($p4.Range.End - 2, $p4.Range.End - 1)
$placeholder.Delete()
